$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.537.63'
$ws.Range('E2').Value = '  +6.53%  '

$ws.Range('D3').Value = '3.323.52'
$ws.Range('E3').Value = '  +2.24%  '

$ws.Range('E4').Value = '  -0.04%  '

$ws.Range('D5').Value = "'409.80"
$ws.Range('E5').Value = '  +3.49%  '

$ws.Range('D6').Value = "'114.82"
$ws.Range('E6').Value = '  +5.69%  '

$ws.Range('D7').Value = '3.316.70'
$ws.Range('E7').Value = '  +2.13%  '

$ws.Range('D8').Value = "'0.572"
$ws.Range('E8').Value = '  -1.30%  '

$ws.Range('E9').Value = '  +0.05%  '

$ws.Range('D10').Value = "'0.624"
$ws.Range('E10').Value = '  -0.19%  '

$ws.Range('E11').Value = '  +19.08%  '

$ws.Range('D12').Value = "'39.93"
$ws.Range('E12').Value = '  +1.70%  '

$ws.Range('E13').Value = '  -0.55%  '

$ws.Range('D14').Value = '3.838.19'
$ws.Range('E14').Value = '  +1.93%  '

$ws.Range('D15').Value = "'8.19"
$ws.Range('E15').Value = '  -1.51%  '

$ws.Range('D16').Value = "'19.10"
$ws.Range('E16').Value = '  -0.05%  '

$ws.Range('D17').Value = '3.307.28'
$ws.Range('E17').Value = '  +1.73%  '

$ws.Range('D18').Value = '60.445.42'
$ws.Range('E18').Value = '  +6.48%  '

$ws.Range('D19').Value = "'0.999"
$ws.Range('E19').Value = '  -3.19%  '

$ws.Range('D20').Value = "'10.77"
$ws.Range('E20').Value = '  +0.50%  '

$ws.Range('E21').Value = '  +5.77%  '

$ws.Range('D22').Value = "'3.36"
$ws.Range('E22').Value = '  -0.01%  '

$ws.Range('D23').Value = "'12.40"
$ws.Range('E23').Value = '  -4.00%  '

$ws.Range('D24').Value = "'295.42"
$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').Value = "'73.95"
$ws.Range('E25').Value = '  -0.38%  '

$ws.Range('D26').Value = "'3.11"
$ws.Range('E26').Value = '  -1.95%  '

$ws.Range('D27').Value = "'29.17"
$ws.Range('E27').Value = '  +3.67%  '

$ws.Range('D28').Value = "'7.62"
$ws.Range('E28').Value = '  +4.00%  '

$ws.Range('E29').Value = '  -2.60%  '

$ws.Range('E30').Value = '  +2.31%  '

$ws.Range('D31').Value = "'7.49"
$ws.Range('E31').Value = '  -1.74%  '

$ws.Range('D32').Value = "'0.114"
$ws.Range('E32').Value = '  +4.44%  '

$ws.Range('E33').Value = '  +0.02%  '

$ws.Range('D34').Value = "'11.24"
$ws.Range('E34').Value = '  -0.20%  '

$ws.Range('B35').Value = 'InjectiveProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D35').Value = "'40.86"
$ws.Range('E35').Value = '  +2.49%  '

$ws.Range('B36').Value = 'Toncoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D36').Value = "'2.47"
$ws.Range('E36').Value = '  +16.06%  '

$ws.Range('D37').Value = "'0.0490"
$ws.Range('E37').Value = '  +0.64%  '

$ws.Range('D38').Value = "'52.16"
$ws.Range('E38').Value = '  +1.56%  '

$ws.Range('D39').Value = "'0.997"
$ws.Range('E39').Value = '  -0.21%  '

$ws.Range('D40').Value = "'3.05"
$ws.Range('E40').Value = '  +5.30%  '

$ws.Range('D41').Value = "'3.36"
$ws.Range('E41').Value = '  -3.59%  '

$ws.Range('D42').Value = "'133.81"
$ws.Range('E42').Value = '  -2.20%  '

$ws.Range('D43').Value = "'0.291"
$ws.Range('E43').Value = '  +3.54%  '

$ws.Range('E44').Value = '  -1.74%  '

$ws.Range('D45').Value = "'1.89"
$ws.Range('E45').Value = '  -0.48%  '

$ws.Range('D46').Value = "'3.83"
$ws.Range('E46').Value = '  -3.77%  '

$ws.Range('D47').Value = "'16.23"
$ws.Range('E47').Value = '  -4.97%  '

$ws.Range('D48').Value = "'2.20"
$ws.Range('E48').Value = '  +2.76%  '

$ws.Range('D49').Value = "'20.99"
$ws.Range('E49').Value = '  -5.90%  '

$ws.Range('D50').Value = '2.134.40'
$ws.Range('E50').Value = '  -0.99%  '

$ws.Range('D51').Value = '3.638.14'
$ws.Range('E51').Value = '  +1.77%  '
